# Listas sem duplicação de professores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E12").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("F21").Value = "-"
